$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 572
$ws.Range("F3").Value = 246
$ws.Range("F4").Value = 25
$ws.Range("F5").Value = 730
$ws.Range("F6").Value = 353
$ws.Range("F8").Value = 141
$ws.Range("F10").Value = 210
$ws.Range("F11").Value = 5846
$ws.Range("F12").Value = 48
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 488
$ws.Range("F16").Value = 541
$ws.Range("F17").Value = 345
$ws.Range("F18").Value = 419
$ws.Range("F22").Value = 65
$ws.Range("F24").Value = 303
$ws.Range("F25").Value = 1006
$ws.Range("F26").Value = 62
$ws.Range("F27").Value = 1776
$ws.Range("F28").Value = 456

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 658
$ws.Range("F5").Value = 263
$ws.Range("F6").Value = 298

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 207

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 207
$ws.Range("F3").Value = 572
$ws.Range("F4").Value = 246
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 730
$ws.Range("F7").Value = 658
$ws.Range("F8").Value = 353
$ws.Range("F10").Value = 141
$ws.Range("F12").Value = 210
$ws.Range("F13").Value = 5846
$ws.Range("F14").Value = 48
$ws.Range("F15").Value = 37
$ws.Range("F17").Value = 488
$ws.Range("F19").Value = 541
$ws.Range("F20").Value = 345
$ws.Range("F21").Value = 419
$ws.Range("F25").Value = 263
$ws.Range("F26").Value = 298
$ws.Range("F32").Value = 65
$ws.Range("F34").Value = 303
$ws.Range("F35").Value = 1006
$ws.Range("F36").Value = 62
$ws.Range("F37").Value = 1776
$ws.Range("F38").Value = 456
